$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "0.3400" or "8.500"
# keep their exact formatting instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.693.17'
$ws.Range("E2").Value = '  +6.97%  '
$ws.Range("D3").Value = '1.735.48'
$ws.Range("E3").Value = '  +3.85%  '
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '332.49'
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").Value = '0.9967'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.3731'
$ws.Range("E7").Value = '  +2.14%  '
$ws.Range("D8").Value = '0.3400'
$ws.Range("E8").Value = '  +4.53%  '
$ws.Range("D9").Value = '48.18'
$ws.Range("E9").Value = '  +1.76%  '
$ws.Range("D10").Value = '1.187'
$ws.Range("E10").Value = '  +3.58%  '
$ws.Range("D11").Value = '0.07465'
$ws.Range("E11").Value = '  +5.39%  '
$ws.Range("D12").Value = '0.9976'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '6.423'
$ws.Range("E13").Value = '  +5.44%  '
$ws.Range("D14").Value = '20.23'
$ws.Range("E14").Value = '  +2.55%  '
$ws.Range("D15").Value = '7.079'
$ws.Range("E15").Value = '  +6.94%  '
$ws.Range("D16").Value = '1.732.77'
$ws.Range("E16").Value = '  +3.81%  '
$ws.Range("D17").Value = '0.00001077'
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("D18").Value = '0.06716'
$ws.Range("E18").Value = '  +1.71%  '
$ws.Range("D19").Value = '82.54'
$ws.Range("E19").Value = '  +4.38%  '
$ws.Range("D20").Value = '0.9969'
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '6.241'
$ws.Range("E21").Value = '  +4.97%  '
$ws.Range("D22").Value = '16.64'
$ws.Range("E22").Value = '  +4.43%  '
$ws.Range("D23").Value = '12.74'
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("D24").Value = '26.675.74'
$ws.Range("E24").Value = '  +7.01%  '
$ws.Range("D25").Value = '2.443'
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("D26").Value = '1.506'
$ws.Range("E26").Value = '  +26.50%  '
$ws.Range("D27").Value = '2.432'
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '151.15'
$ws.Range("E28").Value = '  +1.48%  '
$ws.Range("D29").Value = '19.49'
$ws.Range("E29").Value = '  +3.66%  '
$ws.Range("D30").Value = '1.926.29'
$ws.Range("E30").Value = '  +4.02%  '
$ws.Range("D31").Value = '132.15'
$ws.Range("E31").Value = '  +4.73%  '
$ws.Range("D32").Value = '4.110'
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("D33").Value = '6.045'
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("D34").Value = '0.08632'
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("D35").Value = '1.699'
$ws.Range("E35").Value = '  +3.62%  '
$ws.Range("D36").Value = '12.84'
$ws.Range("E36").Value = '  +4.67%  '
$ws.Range("D37").Value = '5.404'
$ws.Range("E37").Value = '  +4.27%  '
$ws.Range("D38").Value = '0.02346'
$ws.Range("E38").Value = '  +3.61%  '
$ws.Range("D39").Value = '0.2175'
$ws.Range("E39").Value = '  +4.18%  '
$ws.Range("D40").Value = '0.06229'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("D41").Value = '8.500'
$ws.Range("E41").Value = '  +2.86%  '
$ws.Range("D42").Value = '1.224'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '0.6239'
$ws.Range("E43").Value = '  +4.69%  '
$ws.Range("D44").Value = '14.30'
$ws.Range("E44").Value = '  +5.72%  '
$ws.Range("D45").Value = '0.9960'
$ws.Range("E45").Value = '  -0.18%  '
$ws.Range("D46").Value = '3.904'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").Value = '0.6064'
$ws.Range("E47").Value = '  +6.36%  '
$ws.Range("D48").Value = '129.09'
$ws.Range("E48").Value = '  +2.71%  '
$ws.Range("D49").Value = '2.054'
$ws.Range("E49").Value = '  +4.94%  '
$ws.Range("D50").Value = '0.07197'
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("D51").Value = '77.39'
$ws.Range("E51").Value = '  +3.10%  '
